$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the tnrsdate column (Y) for rows 2 through 12 from 2025-09-05 (45905) to 2025-09-09 (45909)
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 25).Value = 45909
}
